$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Películas")
$ws.Activate()

# Insert a new row at 72, shifting the existing rows 72-74 down to 73-75.
$ws.Rows.Item(72).Insert()

# Fill in the new entry: "Estado eléctrico"
$ws.Range("B72").Value = "Estado eléctrico"
$ws.Range("C72").Formula = "=AVERAGE(D72,E72,F72,G72,H72,H72,I72)"
$ws.Range("D72").Value = 6
$ws.Range("E72").Value = 4
$ws.Range("F72").Value = 4
$ws.Range("G72").Value = 6
$ws.Range("H72").Value = 5.9
$ws.Range("I72").Value = 5.0999999999999996

# I72 should render with the "General" number format (like I63), not the
# 0.0 custom format the inserted row inherited from the row above.
$ws.Range("I63").Copy()
$ws.Range("I72").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-apply left alignment on B63 so it collapses back onto the common
# "horizontal left" style shared by the rest of column B.
$ws.Range("B63").HorizontalAlignment = -4131

# Resize the table/ListObject so it covers the new row too.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B2:I75"))

# Update the selection to match where the new row landed.
$ws.Range("C75").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "done"
